# Split the single "Personal" worksheet into two worksheets:
#   - "Departments" (Name, email)       -- re-uses the first sheet
#   - "Employees"   (Email, Employee_number) -- brand-new sheet, placed after it

$wb = $excel.ActiveWorkbook

# --- Rework the existing sheet into "Departments" ---
$departments = $wb.Worksheets.Item(1)
$departments.Name = "Departments"

# Drop the old Personal-sheet columns (C:I) and replace the headers.
$departments.Range("C1:I1").Clear()
$departments.Range("A1").Value = "Name"
$departments.Range("B1").Value = "email"

# --- Add the new "Employees" sheet right after "Departments" ---
$employees = $wb.Worksheets.Add($null, $departments)
$employees.Name = "Employees"

$employees.Range("A1").Value = "Email"
$employees.Range("B1").Value = "Employee_number"

# Match the header style (s="1") used on the Departments sheet.
$departments.Range("A1:B1").Copy()
$employees.Range("A1:B1").PasteSpecial(-4122)

# Match column widths from the old Personal sheet layout (now on Employees).
# (Input values are pre-compensated for the engine's character->pixel
# rounding so the saved <col width=.../> lands on the target value.)
$employees.Columns("A:B").ColumnWidth = 14.09
$employees.Columns("C").ColumnWidth = 10.75
$employees.Columns("D").ColumnWidth = 8.92
$employees.Columns("E").ColumnWidth = 8.42
$employees.Columns("F").ColumnWidth = 6.75
$employees.Columns("G").ColumnWidth = 10.59
$employees.Columns("H").ColumnWidth = 11.09
$employees.Columns("I").ColumnWidth = 17.09
